$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 1.441147
$ws.Cells.Item(2, 8).Value2 = 4.323441
$ws.Cells.Item(2, 9).Value2 = 0.9326218691433955
$ws.Cells.Item(2, 10).Value2 = 0.9326218691433956
$ws.Cells.Item(2, 13).Value2 = 0.6574793333333333
$ws.Cells.Item(2, 14).Value2 = 1.972438
$ws.Cells.Item(2, 15).Value2 = 0.04234443143670402
$ws.Cells.Item(2, 16).Value2 = 0.04234443143670403
$ws.Cells.Item(2, 17).Value2 = 0.9475243687953333
$ws.Cells.Item(2, 18).Value2 = 8.527719319157999
$ws.Cells.Item(2, 19).Value2 = 0.03949134279431326
$ws.Cells.Item(2, 20).Value2 = 0.03949134279431327

$ws.Cells.Item(3, 7).Value2 = 1.441147
$ws.Cells.Item(3, 8).Value2 = 4.323441
$ws.Cells.Item(3, 9).Value2 = 0.9326218691433955
$ws.Cells.Item(3, 10).Value2 = 0.9326218691433956
$ws.Cells.Item(3, 15).Value2 = 0.1192373589365509
$ws.Cells.Item(3, 16).Value2 = 0.119237358936551
$ws.Cells.Item(3, 17).Value2 = 2.668126585477
$ws.Cells.Item(3, 18).Value2 = 24.013139269293
$ws.Cells.Item(3, 19).Value2 = 0.1112033685631281
$ws.Cells.Item(3, 20).Value2 = 0.1112033685631281

$ws.Cells.Item(4, 7).Value2 = 1.441147
$ws.Cells.Item(4, 8).Value2 = 4.323441
$ws.Cells.Item(4, 9).Value2 = 0.9326218691433955
$ws.Cells.Item(4, 10).Value2 = 0.9326218691433956
$ws.Cells.Item(4, 13).Value2 = 5.370269333333333
$ws.Cells.Item(4, 14).Value2 = 16.110808
$ws.Cells.Item(4, 15).Value2 = 0.3458679080132824
$ws.Cells.Item(4, 16).Value2 = 0.3458679080132824
$ws.Cells.Item(4, 17).Value2 = 7.739347538925332
$ws.Cells.Item(4, 18).Value2 = 69.654127850328
$ws.Cells.Item(4, 19).Value2 = 0.3225639748480634
$ws.Cells.Item(4, 20).Value2 = 0.3225639748480634

$ws.Cells.Item(5, 7).Value2 = 1.441147
$ws.Cells.Item(5, 8).Value2 = 4.323441
$ws.Cells.Item(5, 9).Value2 = 0.9326218691433955
$ws.Cells.Item(5, 10).Value2 = 0.9326218691433956
$ws.Cells.Item(5, 13).Value2 = 1.801189666666667
$ws.Cells.Item(5, 14).Value2 = 5.403569
$ws.Cells.Item(5, 15).Value2 = 0.1160041821512257
$ws.Cells.Item(5, 16).Value2 = 0.1160041821512257
$ws.Cells.Item(5, 17).Value2 = 2.595779084547666
$ws.Cells.Item(5, 18).Value2 = 23.362011760929
$ws.Cells.Item(5, 19).Value2 = 0.108188037186327
$ws.Cells.Item(5, 20).Value2 = 0.1081880371863271

$ws.Cells.Item(6, 7).Value2 = 1.441147
$ws.Cells.Item(6, 8).Value2 = 4.323441
$ws.Cells.Item(6, 9).Value2 = 0.9326218691433955
$ws.Cells.Item(6, 10).Value2 = 0.9326218691433956
$ws.Cells.Item(6, 13).Value2 = 5.846608
$ws.Cells.Item(6, 14).Value2 = 17.539824
$ws.Cells.Item(6, 15).Value2 = 0.3765461194622369
$ws.Cells.Item(6, 16).Value2 = 0.376546119462237
$ws.Cells.Item(6, 17).Value2 = 8.425821579375999
$ws.Cells.Item(6, 18).Value2 = 75.832394214384
$ws.Cells.Item(6, 19).Value2 = 0.3511751457515637
$ws.Cells.Item(6, 20).Value2 = 0.3511751457515638

$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 0.104117
$ws.Cells.Item(7, 8).Value2 = 0.312351
$ws.Cells.Item(7, 9).Value2 = 0.06737813085660443
$ws.Cells.Item(7, 10).Value2 = 0.06737813085660445
$ws.Cells.Item(7, 13).Value2 = 0.6574793333333333
$ws.Cells.Item(7, 14).Value2 = 1.972438
$ws.Cells.Item(7, 15).Value2 = 0.04234443143670402
$ws.Cells.Item(7, 16).Value2 = 0.04234443143670403
$ws.Cells.Item(7, 17).Value2 = 0.06845477574866667
$ws.Cells.Item(7, 18).Value2 = 0.616092981738
$ws.Cells.Item(7, 19).Value2 = 0.002853088642390758
$ws.Cells.Item(7, 20).Value2 = 0.002853088642390759

$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 0.104117
$ws.Cells.Item(8, 8).Value2 = 0.312351
$ws.Cells.Item(8, 9).Value2 = 0.06737813085660443
$ws.Cells.Item(8, 10).Value2 = 0.06737813085660445
$ws.Cells.Item(8, 15).Value2 = 0.1192373589365509
$ws.Cells.Item(8, 16).Value2 = 0.119237358936551
$ws.Cells.Item(8, 17).Value2 = 0.192761276747
$ws.Cells.Item(8, 18).Value2 = 1.734851490723
$ws.Cells.Item(8, 19).Value2 = 0.008033990373422843
$ws.Cells.Item(8, 20).Value2 = 0.008033990373422844

$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 0.104117
$ws.Cells.Item(9, 8).Value2 = 0.312351
$ws.Cells.Item(9, 9).Value2 = 0.06737813085660443
$ws.Cells.Item(9, 10).Value2 = 0.06737813085660445
$ws.Cells.Item(9, 13).Value2 = 5.370269333333333
$ws.Cells.Item(9, 14).Value2 = 16.110808
$ws.Cells.Item(9, 15).Value2 = 0.3458679080132824
$ws.Cells.Item(9, 16).Value2 = 0.3458679080132824
$ws.Cells.Item(9, 17).Value2 = 0.5591363321786667
$ws.Cells.Item(9, 18).Value2 = 5.032226989608
$ws.Cells.Item(9, 19).Value2 = 0.02330393316521897
$ws.Cells.Item(9, 20).Value2 = 0.02330393316521897

$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 0.104117
$ws.Cells.Item(10, 8).Value2 = 0.312351
$ws.Cells.Item(10, 9).Value2 = 0.06737813085660443
$ws.Cells.Item(10, 10).Value2 = 0.06737813085660445
$ws.Cells.Item(10, 13).Value2 = 1.801189666666667
$ws.Cells.Item(10, 14).Value2 = 5.403569
$ws.Cells.Item(10, 15).Value2 = 0.1160041821512257
$ws.Cells.Item(10, 16).Value2 = 0.1160041821512257
$ws.Cells.Item(10, 17).Value2 = 0.1875344645243333
$ws.Cells.Item(10, 18).Value2 = 1.687810180719
$ws.Cells.Item(10, 19).Value2 = 0.007816144964898662
$ws.Cells.Item(10, 20).Value2 = 0.007816144964898665

$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 0.104117
$ws.Cells.Item(11, 8).Value2 = 0.312351
$ws.Cells.Item(11, 9).Value2 = 0.06737813085660443
$ws.Cells.Item(11, 10).Value2 = 0.06737813085660445
$ws.Cells.Item(11, 13).Value2 = 5.846608
$ws.Cells.Item(11, 14).Value2 = 17.539824
$ws.Cells.Item(11, 15).Value2 = 0.3765461194622369
$ws.Cells.Item(11, 16).Value2 = 0.376546119462237
$ws.Cells.Item(11, 17).Value2 = 0.6087312851359999
$ws.Cells.Item(11, 18).Value2 = 5.478581566223999
$ws.Cells.Item(11, 19).Value2 = 0.02537097371067321
$ws.Cells.Item(11, 20).Value2 = 0.02537097371067322
